$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.584159851074219
$ws.Range("B1").Value = 4.958879947662354
$ws.Range("C1").Value = 4.926409244537354
$ws.Range("D1").Value = 8.028890609741211
$ws.Range("E1").Value = 3.429958343505859
